$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 6
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = 7
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 0.383
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 12.7

$ws = $wb.Worksheets.Item(2)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 8
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = 9
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 0.096
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 2.4

$ws = $wb.Worksheets.Item(3)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 13
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = 14
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 0.177
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 12.1

$ws = $wb.Worksheets.Item(4)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 8
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = 9
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 0.289
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 16.1

$ws = $wb.Worksheets.Item(5)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 16
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = 17
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 0.182
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 13.4

$ws = $wb.Worksheets.Item(6)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 10
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = 11
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 0.213
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 17.1

$ws = $wb.Worksheets.Item(7)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 11
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = 12
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 0.363
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 19.9

$ws = $wb.Worksheets.Item(8)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 7
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = 8
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 2.557
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 37.4

$ws = $wb.Worksheets.Item(9)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 21
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = 22
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 3.598
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 52.2

$ws = $wb.Worksheets.Item(10)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 32
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = 33
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 19.823
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 78.1

$ws = $wb.Worksheets.Item(11)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 46
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = "N/A"
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 1806.026
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 355.4

$ws = $wb.Worksheets.Item(12)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 29
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = 30
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 3.919
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 71.2

$ws = $wb.Worksheets.Item(13)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 33
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = "N/A"
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 1806.108
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 2078.9

$ws = $wb.Worksheets.Item(14)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 33
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = "N/A"
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 1806.013
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 2075.8

$ws = $wb.Worksheets.Item(15)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 89
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = "N/A"
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 1805.881
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 1188

$ws = $wb.Worksheets.Item(16)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 67
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = "N/A"
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 1806.029
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 1444.8

$ws = $wb.Worksheets.Item(17)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 136
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = "N/A"
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 1806.722
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 1257.4

$ws = $wb.Worksheets.Item(18)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 60
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = "N/A"
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 1806.459
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 2157.1

$ws = $wb.Worksheets.Item(19)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 138
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = "N/A"
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 1806.561
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 1703.8

$ws = $wb.Worksheets.Item(20)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 75
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = "N/A"
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 1805.728
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 2545.7

$ws = $wb.Worksheets.Item(21)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 137
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = "N/A"
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 1805.805
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 2150.7

$ws = $wb.Worksheets.Item(22)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 170
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = "N/A"
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 1805.788
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 2531.6

$ws = $wb.Worksheets.Item(23)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 111
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = "N/A"
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 1805.679
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 2508.9

$ws = $wb.Worksheets.Item(24)
$ws.Range("F1").Value = "Max width SAT"
$ws.Range("G1").Value = 73
$ws.Range("F2").Value = "Min width UNSAT"
$ws.Range("G2").Value = "N/A"
$ws.Range("F3").Value = "Total real time (s)"
$ws.Range("G3").Value = 1805.695
$ws.Range("F4").Value = "Total memory consumed (MB)"
$ws.Range("G4").Value = 3634.8
